$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.163.68'
$ws.Range('E2').Value = '  +8.08%  '
$ws.Range('D3').Value = '1.876.52'
$ws.Range('E3').Value = '  +5.44%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9994'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '250.02'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9989'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4983'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.17%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2854'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +7.60%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06616'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.15%  '
$ws.Range('D10').Value = '1.871.57'
$ws.Range('E10').Value = '  +5.28%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '17.12'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.87%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07214'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.33%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6654'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +8.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '85.58'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +8.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.832'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.75%  '
$ws.Range('D16').Value = '30.176.76'
$ws.Range('E16').Value = '  +8.28%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9987'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.97'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +10.35%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007575'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.28%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9974'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.09%  '
$ws.Range('D21').Value = '2.109.77'
$ws.Range('E21').Value = '  +5.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.798'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.83%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.089'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.31%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.524'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +6.68%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.03'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.05%  '
$ws.Range('B26').Value = 'BitcoinCash'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '136.66'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +24.80%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.81'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +8.15%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.960'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.70%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.400'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.53%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.252'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08642'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.913'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05097'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +8.60%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.133'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +7.95%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6902'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +9.37%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9978'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('E37').Value = '  +2.92%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.342'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +15.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.745'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9623'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.78%  '
$ws.Range('E41').Value = '  +7.39%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.124'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '103.82'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.86%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9987'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4216'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.70%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.518'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +7.25%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1262'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +6.04%  '
$ws.Range('E48').Value = '  +4.43%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '32.71'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.99%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.327'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3746'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +8.99%  '
